# ---------------------------------------------------------------------------
# Adds 5 reference list sheets (Ref_TypesDocument, Ref_Thematiques,
# Ref_MotsCles, Ref_SourcesDocument, Ref_DomainesMetier) and 5 new columns
# (Type_Document, Domaine_Metier, Source_Document, Thematiques, Mots_Cles) to
# the Validation_Questions sheet, with data-validation dropdown lists wired
# up to the new reference sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Add-RefSheet($name, $values) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $sheet = $wb.Worksheets.Add($null, $lastSheet)
    $sheet.Name = $name
    $r = 1
    foreach ($v in $values) {
        $sheet.Cells.Item($r, 1).Value = $v
        $r = $r + 1
    }
    return $sheet
}

# --- 1. New reference sheets (in the order they appear in the workbook) ---

$typesDocument = @(
    "Type de document",
    "Actualités",
    "Assurances",
    "Convention collectives Notariat",
    "Directives CSN",
    "Lois et règlements"
)
Add-RefSheet "Ref_TypesDocument" $typesDocument | Out-Null

$thematiques = @(
    "Thématique",
    "acte authentique",
    "avenant",
    "circulaire",
    "congés payés",
    "conseil supérieur du notariat",
    "convention collective",
    "cybersécurité",
    "formation professionnelle",
    "harcèlement",
    "intéressement",
    "lcb-ft",
    "licenciement",
    "minute",
    "opco",
    "prévoyance",
    "rgpd",
    "rémunération",
    "tarification",
    "égalité professionnelle"
)
Add-RefSheet "Ref_Thematiques" $thematiques | Out-Null

$motsCles = @(
    "Mot-clé",
    "CRIDON",
    "FAQ notariale",
    "actualité juridique",
    "administration",
    "assurance professionnelle",
    "base de connaissances",
    "bonnes pratiques",
    "carrière notaire",
    "circulaire CSN",
    "complémentaire santé",
    "conformité",
    "congés payés",
    "consultation juridique",
    "convention collective",
    "cyber-risques",
    "cybersécurité",
    "doctrine",
    "documentation",
    "documentation métier",
    "droit social",
    "déclarations",
    "expertise notariale",
    "formation",
    "formation professionnelle",
    "harcèlement au travail",
    "instructions professionnelles",
    "législation",
    "négociation collective",
    "obligations fiscales",
    "parcours professionnel",
    "participation aux bénéfices",
    "procédure disciplinaire",
    "protection des données",
    "protection données",
    "prévention harcèlement",
    "prévoyance",
    "questions-réponses",
    "recherche juridique",
    "responsabilité civile",
    "réglementation notariale",
    "rémunération",
    "textes réglementaires",
    "transactions immobilières",
    "veille professionnelle",
    "égalité professionnelle"
)
Add-RefSheet "Ref_MotsCles" $motsCles | Out-Null

$sourcesDocument = @(
    "Source de document",
    "accord_branche",
    "assurance",
    "avenant_ccn",
    "circulaire_csn",
    "conformite",
    "decret_ordonnance",
    "fil_info",
    "guide_pratique"
)
Add-RefSheet "Ref_SourcesDocument" $sourcesDocument | Out-Null

$domainesMetier = @(
    "Domaine métier",
    "ASSURANCES",
    "DEONTOLOGIE",
    "RH"
)
Add-RefSheet "Ref_DomainesMetier" $domainesMetier | Out-Null

# --- 2. Insert 5 new columns into Validation_Questions (before "Difficulte") ---

$ws = $wb.Worksheets.Item("Validation_Questions")
$ws.Columns("D:H").Insert()

# The insert carries the neighbouring row style onto the new stub cells;
# wipe them out entirely so the new columns start fully blank.
$ws.Range("D2:H21").Clear()

# --- 3. Header labels for the 5 new columns ---

$ws.Range("D1:H1").Style = "Normal"
$ws.Range("D1").Value = "Type_Document"
$ws.Range("E1").Value = "Domaine_Metier"
$ws.Range("F1").Value = "Source_Document"
$ws.Range("G1").Value = "Thematiques"
$ws.Range("H1").Value = "Mots_Cles"

# --- 4. Data validation dropdowns for the new columns ---

$dv = $ws.Range("D2:D21").Validation
$dv.Add(3, 1, 1, "=Ref_TypesDocument!`$A`$2:`$A`$6")
$dv.InputTitle = "Type de document"
$dv.InputMessage = "Choisissez un type de document"
$dv.IgnoreBlank = $true
$dv.ShowInput = $true
$dv.ShowError = $false

$dv = $ws.Range("E2:E21").Validation
$dv.Add(3, 1, 1, "=Ref_DomainesMetier!`$A`$2:`$A`$4")
$dv.InputTitle = "Domaine métier"
$dv.InputMessage = "Choisissez un domaine métier"
$dv.IgnoreBlank = $true
$dv.ShowInput = $true
$dv.ShowError = $false

$dv = $ws.Range("F2:F21").Validation
$dv.Add(3, 1, 1, "=Ref_SourcesDocument!`$A`$2:`$A`$9")
$dv.InputTitle = "Source de document"
$dv.InputMessage = "Choisissez une source de document"
$dv.IgnoreBlank = $true
$dv.ShowInput = $true
$dv.ShowError = $false

$dv = $ws.Range("G2:G21").Validation
$dv.Add(3, 1, 1, "=Ref_Thematiques!`$A`$2:`$A`$20")
$dv.InputTitle = "Thématiques"
$dv.InputMessage = "Choisissez une ou plusieurs thématiques (séparées par des virgules)"
$dv.IgnoreBlank = $true
$dv.ShowInput = $true
$dv.ShowError = $false

$dv = $ws.Range("H2:H21").Validation
$dv.Add(3, 1, 1, "=Ref_MotsCles!`$A`$2:`$A`$46")
$dv.InputTitle = "Mots-clés"
$dv.InputMessage = "Choisissez un ou plusieurs mots-clés (séparés par des virgules)"
$dv.IgnoreBlank = $true
$dv.ShowInput = $true
$dv.ShowError = $false

Write-Output "edit applied"
